$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 154, shifting existing rows 154:193 down to 155:194.
$ws.Rows("154").Insert()

# Populate the newly inserted row 154 with its data (matches the
# surrounding rows' constant columns, with its own date/price figures).
$ws.Range("A154").Value = 8
$ws.Range("B154").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C154").Value = 'Coquimbo'
$ws.Range("D154").Value = 44943
$ws.Range("E154").Value = 4
$ws.Range("F154").Value = 100112040
$ws.Range("G154").Value = 'Cilantro'
$ws.Range("H154").Value = 'Sin especificar'
$ws.Range("I154").Value = 'Primera'
$ws.Range("J154").Value = 2000
$ws.Range("K154").Value = 3000
$ws.Range("L154").Value = 3500
$ws.Range("M154").Value = 3250
$ws.Range("N154").Value = '$/atado 1 a 1,5 kilos'
$ws.Range("O154").Value = 'Provincia del Elqu' + [char]0xED
$ws.Range("P154").Value = 2167
$ws.Range("Q154").Value = 1.5
$ws.Range("R154").Value = 'Hortaliza'

# Match the date-format style used by the rest of column D.
$ws.Range("D154").NumberFormat = $ws.Range("D155").NumberFormat
